$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 18735
$ws.Range("I100").Value = 21682
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 21682
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -21141
$ws.Range("N100").Value = -5082
$ws.Range("H120").Value = 24332.334
$ws.Range("J120").Value = 24332.334
$ws.Range("L120").Value = 24332.334
$ws.Range("N120").Value = -34008.334

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 23000
$ws.Range("J51").Value = 23000
$ws.Range("L51").Value = 23000
$ws.Range("N51").Value = -24512
$ws.Range("H63").Value = 2426.158
$ws.Range("I63").Value = 1991.1666
$ws.Range("J63").Value = 2626.923
$ws.Range("K63").Value = 1991.1666
$ws.Range("L63").Value = 2626.923
$ws.Range("M63").Value = -1305.1666
$ws.Range("N63").Value = -3998.923
$ws.Range("H66").Value = 2426.158
$ws.Range("I66").Value = 1991.1666
$ws.Range("J66").Value = 2626.923
$ws.Range("K66").Value = 9955.833000000001
$ws.Range("L66").Value = 13134.615
$ws.Range("M66").Value = -6523.833000000001
$ws.Range("N66").Value = -19998.615
$ws.Range("H97").Value = 1195.0555
$ws.Range("I97").Value = 1150
$ws.Range("J97").Value = 1352.75
$ws.Range("K97").Value = 1150
$ws.Range("L97").Value = 1352.75
$ws.Range("M97").Value = -654
$ws.Range("N97").Value = -2344.75
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 1133260.6
$ws.Range("I132").Value = 1935.973
$ws.Range("J132").Value = 3923861.5
$ws.Range("K132").Value = 5807.919
$ws.Range("L132").Value = 11771584.5
$ws.Range("M132").Value = -3277.919
$ws.Range("N132").Value = -11776644.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = -6755
$ws.Range("H105").Value = 996.74225
$ws.Range("I105").Value = 996.0273999999999
$ws.Range("J105").Value = 998.9167
$ws.Range("K105").Value = 996.0273999999999
$ws.Range("L105").Value = 998.9167
$ws.Range("M105").Value = 750.9726000000001
$ws.Range("N105").Value = -4492.9167
$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 45000
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 35000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 35000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 35000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -35348
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 200
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -26
$ws.Range("N25").ClearContents()
$ws.Range("H99").Value = 45461784
$ws.Range("I99").Value = 250026430
$ws.Range("J99").Value = 2975.2222
$ws.Range("K99").Value = 250026430
$ws.Range("L99").Value = 2975.2222
$ws.Range("M99").Value = -250024932
$ws.Range("N99").Value = -5971.2222
$ws.Range("H126").Value = 45461784
$ws.Range("I126").Value = 250026430
$ws.Range("J126").Value = 2975.2222
$ws.Range("K126").Value = 750079290
$ws.Range("L126").Value = 8925.6666
$ws.Range("M126").Value = -750076820
$ws.Range("N126").Value = -13865.6666
$ws.Range("H132").Value = 5557076.5
$ws.Range("I132").Value = 1184.075
$ws.Range("J132").Value = 16668862
$ws.Range("K132").Value = 3552.225
$ws.Range("L132").Value = 50006586
$ws.Range("M132").Value = -1022.225
$ws.Range("N132").Value = -50011646
$ws.Range("H134").Value = 1254.2667
$ws.Range("I134").Value = 945.12
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 2835.36
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -300.3600000000001
$ws.Range("N134").Value = -13470

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2000
$ws.Range("I87").Value = 2000
$ws.Range("K87").Value = 6000
$ws.Range("M87").Value = -4752
$ws.Range("H90").Value = 2000
$ws.Range("I90").Value = 2000
$ws.Range("K90").Value = 18000
$ws.Range("M90").Value = -11760
$ws.Range("H99").Value = 1504.1666
$ws.Range("I99").Value = 1005
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 3015
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -769
$ws.Range("N99").Value = -16492
$ws.Range("H104").Value = 2159.8
$ws.Range("I104").Value = 899.5
$ws.Range("J104").Value = 3000
$ws.Range("K104").Value = 2698.5
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = -77.5
$ws.Range("N104").Value = -14242
$ws.Range("H109").Value = 3564.5334
$ws.Range("I109").Value = 782.2857
$ws.Range("J109").Value = 5999
$ws.Range("K109").Value = 2346.8571
$ws.Range("L109").Value = 17997
$ws.Range("M109").Value = -1306.8571
$ws.Range("N109").Value = -20077
$ws.Range("H126").Value = 31253312
$ws.Range("I126").Value = 500000000
$ws.Range("J126").Value = 3533.3333
$ws.Range("K126").Value = 1500000000
$ws.Range("L126").Value = 10599.9999
$ws.Range("M126").Value = -1499995060
$ws.Range("N126").Value = -20479.9999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3250
$ws.Range("I97").Value = 3080
$ws.Range("J97").Value = 3371.4285
$ws.Range("K97").Value = 3080
$ws.Range("L97").Value = 3371.4285
$ws.Range("M97").Value = -2584
$ws.Range("N97").Value = -4363.4285
$ws.Range("H132").Value = 6152.1924
$ws.Range("I132").Value = 1619.8889
$ws.Range("J132").Value = 16349.875
$ws.Range("K132").Value = 4859.6667
$ws.Range("L132").Value = 49049.625
$ws.Range("M132").Value = -2329.6667
$ws.Range("N132").Value = -54109.625

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2977645.5
$ws.Range("I46").Value = 3788723.2
$ws.Range("J46").Value = 3693.3333
$ws.Range("K46").Value = 3788723.2
$ws.Range("L46").Value = 3693.3333
$ws.Range("M46").Value = -3788535.2
$ws.Range("N46").Value = -4069.3333
$ws.Range("H93").Value = 1062.6957
$ws.Range("I93").Value = 1086.8334
$ws.Range("J93").Value = 1036.3636
$ws.Range("K93").Value = 1086.8334
$ws.Range("L93").Value = 1036.3636
$ws.Range("M93").Value = 161.1666
$ws.Range("N93").Value = -3532.3636
$ws.Range("H111").Value = 59800
$ws.Range("J111").Value = 59800
$ws.Range("L111").Value = 59800
$ws.Range("N111").Value = -67980
$ws.Range("H122").Value = 18167.334
$ws.Range("I122").Value = 18167.334
$ws.Range("K122").Value = 54502.00199999999
$ws.Range("M122").Value = -52052.00199999999
$ws.Range("H132").Value = 27219294
$ws.Range("I132").Value = 54424600
$ws.Range("J132").Value = 13986.762
$ws.Range("K132").Value = 163273800
$ws.Range("L132").Value = 41960.286
$ws.Range("M132").Value = -163271270
$ws.Range("N132").Value = -47020.286

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2136
$ws.Range("I96").Value = 1624.5
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 1624.5
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = -251.5
$ws.Range("N96").Value = -6246
